$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Collapse the split "P_DEMO.htm" hyperlink text back into one run,
#    keeping the Hyperlink character style.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute(
    "https://wwwn.cdc.gov/Nchs/Nhanes/2017-2018/P_D" + "E" + "MO.htm",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://wwwn.cdc.gov/Nchs/Nhanes/2017-2018/P_DEMO.htm", 2)

$rng1b = $d.Content
$rng1b.Find.Execute(
    "https://wwwn.cdc.gov/Nchs/Nhanes/2017-2018/P_DEMO.htm",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng1b.Find.Found) {
    $rng1b.Style = "Hyperlink"
}

# ---------------------------------------------------------------------
# 2) Collapse the split "P_DIQ.htm" hyperlink text back into one run,
#    keeping the Hyperlink character style.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(
    "https://wwwn.cdc.gov/Nchs/Nhanes/2017-20" + "1" + "8/P_DIQ.htm",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://wwwn.cdc.gov/Nchs/Nhanes/2017-2018/P_DIQ.htm", 2)

$rng2b = $d.Content
$rng2b.Find.Execute(
    "https://wwwn.cdc.gov/Nchs/Nhanes/2017-2018/P_DIQ.htm",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2b.Find.Found) {
    $rng2b.Style = "Hyperlink"
}

# ---------------------------------------------------------------------
# 3) Split the paragraph that starts right after the page break: the new
#    leading paragraph keeps the page-break marker plus a new sentence
#    about the temporary GitHub access token; the original sentence
#    ("Sometimes the missing values...") becomes its own paragraph with
#    no page-break marker, its run structure (two runs) preserved.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute(
    "Sometimes the missing values in the data include inapplicable cases (such as A1C reported as missing for non-diabetics). Point is, the data has to make sense.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng3.Find.Found) {
    $xmlPkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
        "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
        "<pkg:xmlData>" +
        "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:body>" +
        "<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space='preserve'>Using temporary access token (only available till August 22, 2024) for GitHub access, remote access denied otherwise. Remote access token: </w:t></w:r></w:p>" +
        "<w:p><w:r><w:t>Sometimes the missing values in the data include inapplicable cases (such as A1C reported as missing for non-diabetics). Point is, the data has to make sense.</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> When we ask if diabetics are all diagnosed with the same set of questions regardless of demographics, we need to make sure the information isn" + [char]8217 + "t missing or unaccounted simply because they" + [char]8217 + "re non-diabetics.</w:t></w:r></w:p>" +
        "</w:body>" +
        "</w:document>" +
        "</pkg:xmlData></pkg:part></pkg:package>"
    $rng3.InsertXML($xmlPkg)
}
